# Auto-generated edit script for Balmung_Profits workbook update
# Applies per-cell numeric updates across 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (69 changes) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 319
$ws.Range("H6").Value = 498.41666
$ws.Range("I6").Value = 488.2
$ws.Range("J6").Value = 549.5
$ws.Range("K6").Value = 1464.6
$ws.Range("L6").Value = 1648.5
$ws.Range("M6").Value = -1352.6
$ws.Range("N6").Value = -1872.5
$ws.Range("H18").Value = 3774.1667
$ws.Range("I18").Value = 3774.1667
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 3774.1667
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -3490.1667
$ws.Range("N18").ClearContents()
$ws.Range("H32").Value = 27779644
$ws.Range("I32").Value = 1505.6
$ws.Range("K32").Value = 1505.6
$ws.Range("M32").Value = -1179.6
$ws.Range("H33").Value = 8208812
$ws.Range("I33").Value = 10943296
$ws.Range("K33").Value = 10943296
$ws.Range("M33").Value = -10943067
$ws.Range("H39").Value = 1823.0769
$ws.Range("J39").Value = 2700
$ws.Range("L39").Value = 8100
$ws.Range("N39").Value = -8692
$ws.Range("H40").Value = 1279.375
$ws.Range("J40").Value = 1298.0667
$ws.Range("L40").Value = 1298.0667
$ws.Range("N40").Value = -1648.0667
$ws.Range("H53").Value = 66667070
$ws.Range("I53").Value = 187.16667
$ws.Range("K53").Value = 187.16667
$ws.Range("M53").Value = 449.83333
$ws.Range("H64").Value = 5723.5
$ws.Range("I64").Value = 4057.8
$ws.Range("K64").Value = 4057.8
$ws.Range("M64").Value = -3809.8
$ws.Range("H67").Value = 5723.5
$ws.Range("I67").Value = 4057.8
$ws.Range("K67").Value = 4057.8
$ws.Range("M67").Value = -3199.8
$ws.Range("H74").Value = 4616.0713
$ws.Range("I74").Value = 4616.0713
$ws.Range("K74").Value = 4616.0713
$ws.Range("M74").Value = -3680.0713
$ws.Range("H77").Value = 4616.0713
$ws.Range("I77").Value = 4616.0713
$ws.Range("K77").Value = 23080.3565
$ws.Range("M77").Value = -18400.3565
$ws.Range("H98").Value = 3365.5588
$ws.Range("I98").Value = 1974.0952
$ws.Range("J98").Value = 5613.3076
$ws.Range("K98").Value = 1974.0952
$ws.Range("L98").Value = 5613.3076
$ws.Range("M98").Value = -476.0952
$ws.Range("N98").Value = -8609.3076
$ws.Range("H122").Value = 3365.5588
$ws.Range("I122").Value = 1974.0952
$ws.Range("J122").Value = 5613.3076
$ws.Range("K122").Value = 5922.2856
$ws.Range("L122").Value = 16839.9228
$ws.Range("M122").Value = -3472.2856
$ws.Range("N122").Value = -21739.9228
$ws.Range("H138").Value = 5074.15
$ws.Range("J138").Value = 3561.6785
$ws.Range("L138").Value = 10685.0355
$ws.Range("N138").Value = -20965.0355

# ---- Sheet: ARM (15 changes) ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 454
$ws.Range("I5").Value = 459
$ws.Range("K5").Value = 459
$ws.Range("M5").Value = -347
$ws.Range("H122").Value = 1069.091
$ws.Range("I122").Value = 1094.4
$ws.Range("J122").Value = 816
$ws.Range("K122").Value = 3283.2
$ws.Range("L122").Value = 2448
$ws.Range("M122").Value = -833.2000000000003
$ws.Range("N122").Value = -7348
$ws.Range("H132").Value = 2238.524
$ws.Range("I132").Value = 2034.0513
$ws.Range("K132").Value = 6102.1539
$ws.Range("M132").Value = -3572.1539

# ---- Sheet: BSM (20 changes) ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 454
$ws.Range("I4").Value = 459
$ws.Range("K4").Value = 459
$ws.Range("M4").Value = -344
$ws.Range("H7").Value = 458.8
$ws.Range("J7").Value = 1100
$ws.Range("L7").Value = 1100
$ws.Range("N7").Value = -1326
$ws.Range("H22").Value = 1110
$ws.Range("I22").Value = 840
$ws.Range("K22").Value = 840
$ws.Range("M22").Value = -667
$ws.Range("H122").Value = 49835.082
$ws.Range("J122").Value = 49835.082
$ws.Range("L122").Value = 49835.082
$ws.Range("N122").Value = -59635.082
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

# ---- Sheet: CRP (34 changes) ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 3990.5
$ws.Range("I10").Value = 2033.5
$ws.Range("K10").Value = 2033.5
$ws.Range("M10").Value = -1894.5
$ws.Range("H22").Value = 1025.7037
$ws.Range("J22").Value = 1256.5
$ws.Range("L22").Value = 1256.5
$ws.Range("N22").Value = -1956.5
$ws.Range("H99").Value = 6668998.5
$ws.Range("J99").Value = 3497.5
$ws.Range("L99").Value = 3497.5
$ws.Range("N99").Value = -6493.5
$ws.Range("H105").Value = 3073.6667
$ws.Range("I105").Value = 2480.6667
$ws.Range("J105").Value = 3666.6667
$ws.Range("K105").Value = 2480.6667
$ws.Range("L105").Value = 3666.6667
$ws.Range("M105").Value = -733.6667000000002
$ws.Range("N105").Value = -7160.6667
$ws.Range("H126").Value = 6668998.5
$ws.Range("J126").Value = 3497.5
$ws.Range("L126").Value = 10492.5
$ws.Range("N126").Value = -15432.5
$ws.Range("H132").Value = 29332.107
$ws.Range("I132").Value = 42201.48
$ws.Range("K132").Value = 126604.44
$ws.Range("M132").Value = -124074.44
$ws.Range("H134").Value = 1622.44
$ws.Range("I134").Value = 1464.2778
$ws.Range("J134").Value = 2029.1428
$ws.Range("K134").Value = 4392.8334
$ws.Range("L134").Value = 6087.428400000001
$ws.Range("M134").Value = -1857.8334
$ws.Range("N134").Value = -11157.4284

# ---- Sheet: CUL (15 changes) ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 3703.5264
$ws.Range("I129").Value = 1217.25
$ws.Range("J129").Value = 5511.727
$ws.Range("K129").Value = 3651.75
$ws.Range("L129").Value = 16535.181
$ws.Range("M129").Value = 1348.25
$ws.Range("N129").Value = -26535.181
$ws.Range("H140").Value = 22224652
$ws.Range("J140").Value = 2400
$ws.Range("L140").Value = 7200
$ws.Range("N140").Value = -17560
$ws.Range("H141").Value = 4854
$ws.Range("I141").Value = 4829.6665
$ws.Range("K141").Value = 14488.9995
$ws.Range("M141").Value = -9308.999500000002

# ---- Sheet: GSM (12 changes) ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2914.1428
$ws.Range("I102").Value = 2444.3
$ws.Range("K102").Value = 2444.3
$ws.Range("M102").Value = -822.3000000000002
$ws.Range("H134").Value = 78000.60000000001
$ws.Range("J134").Value = 78000.60000000001
$ws.Range("L134").Value = 234001.8
$ws.Range("N134").Value = -239071.8
$ws.Range("H136").Value = 36500
$ws.Range("J136").Value = 36500
$ws.Range("L136").Value = 109500
$ws.Range("N136").Value = -114600

# ---- Sheet: LTW (35 changes) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1049.5
$ws.Range("J55").Value = 1343.3846
$ws.Range("L55").Value = 1343.3846
$ws.Range("N55").Value = -1689.3846
$ws.Range("H68").Value = 2440.5
$ws.Range("I68").Value = 1587.3334
$ws.Range("K68").Value = 1587.3334
$ws.Range("M68").Value = -838.3334
$ws.Range("H71").Value = 2440.5
$ws.Range("I71").Value = 1587.3334
$ws.Range("K71").Value = 7936.666999999999
$ws.Range("M71").Value = -4192.666999999999
$ws.Range("H93").Value = 1990.6
$ws.Range("I93").Value = 1738.375
$ws.Range("K93").Value = 1738.375
$ws.Range("M93").Value = -490.375
$ws.Range("H104").Value = 48290
$ws.Range("J104").Value = 48290
$ws.Range("L104").Value = 48290
$ws.Range("N104").Value = -55278
$ws.Range("H122").Value = 3274.7292
$ws.Range("I122").Value = 2835.1875
$ws.Range("K122").Value = 8505.5625
$ws.Range("M122").Value = -6055.5625
$ws.Range("H132").Value = 3038.8823
$ws.Range("I132").Value = 2622.4092
$ws.Range("J132").Value = 3802.4167
$ws.Range("K132").Value = 7867.2276
$ws.Range("L132").Value = 11407.2501
$ws.Range("M132").Value = -5337.2276
$ws.Range("N132").Value = -16467.2501
$ws.Range("H140").Value = 99999.5
$ws.Range("J140").Value = 99999.5
$ws.Range("L140").Value = 99999.5
$ws.Range("N140").Value = -110359.5

# ---- Sheet: WVR (27 changes) ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 495
$ws.Range("J22").Value = 495
$ws.Range("L22").Value = 495
$ws.Range("N22").Value = -1081
$ws.Range("H62").Value = 2746.5
$ws.Range("I62").Value = 2746.5
$ws.Range("K62").Value = 2746.5
$ws.Range("M62").Value = -2122.5
$ws.Range("H65").Value = 2746.5
$ws.Range("I65").Value = 2746.5
$ws.Range("K65").Value = 13732.5
$ws.Range("M65").Value = -10612.5
$ws.Range("H81").Value = 178629.36
$ws.Range("I81").Value = 2769.3333
$ws.Range("K81").Value = 5538.6666
$ws.Range("M81").Value = -4477.6666
$ws.Range("H84").Value = 178629.36
$ws.Range("I84").Value = 2769.3333
$ws.Range("K84").Value = 27693.333
$ws.Range("M84").Value = -22389.333
$ws.Range("H107").Value = 752427.9399999999
$ws.Range("I107").Value = 565.125
$ws.Range("J107").Value = 4762363
$ws.Range("K107").Value = 1695.375
$ws.Range("L107").Value = 14287089
$ws.Range("M107").Value = 224.625
$ws.Range("N107").Value = -14290929
